$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Template cell used as the source of formatting for new optional (yellow) header cells
$formatSource = $ws.Range("C15")

$cell = $ws.Range('AG15')
$cell.Value = 'biotic_relationship'
$formatSource.Copy()
$cell.PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$cell.AddComment('Free-living or from host (define relationship)') | Out-Null

$cell = $ws.Range('AH15')
$cell.Value = 'dew_point'
$formatSource.Copy()
$cell.PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$cell.AddComment('temperature to which a given parcel of humid air must be cooled, at constant barometric pressure, for water vapor to condense into water.') | Out-Null

$cell = $ws.Range('AI15')
$cell.Value = 'encoded_traits'
$formatSource.Copy()
$cell.PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$cell.AddComment('Traits like antibiotic resistance/xenobiotic degration phenotypes/converting phage genes') | Out-Null

$cell = $ws.Range('AJ15')
$cell.Value = 'estimated_size'
$formatSource.Copy()
$cell.PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$cell.AddComment('Estimated size of genome') | Out-Null

$cell = $ws.Range('AK15')
$cell.Value = 'health_state'
$formatSource.Copy()
$cell.PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$cell.AddComment('Health or disease status of sample at time of collection') | Out-Null

$cell = $ws.Range('AL15')
$cell.Value = 'host'
$formatSource.Copy()
$cell.PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$cell.AddComment('The natural (as opposed to laboratory) host to the organism from which the sample was obtained. Use the full taxonomic name, eg, "Homo sapiens".') | Out-Null

$cell = $ws.Range('AM15')
$cell.Value = 'host_taxid'
$formatSource.Copy()
$cell.PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$cell.AddComment('NCBI taxonomy ID of the host, e.g. 9606') | Out-Null

$cell = $ws.Range('AN15')
$cell.Value = 'indoor_surf'
$formatSource.Copy()
$cell.PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$cell.AddComment('type of indoor surface') | Out-Null

$cell = $ws.Range('AO15')
$cell.Value = 'isolation_source'
$formatSource.Copy()
$cell.PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$cell.AddComment('Describes the physical, environmental and/or local geographical source of the biological sample from which the sample was derived.') | Out-Null

$cell = $ws.Range('AP15')
$cell.Value = 'num_replicons'
$formatSource.Copy()
$cell.PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$cell.AddComment('Reports the number of replicons in a nuclear genome of eukaryotes, in the genome of a bacterium or archaea or the number of segments in a segmented virus. Always applied to the haploid chromosome count of a eukaryote') | Out-Null

$cell = $ws.Range('AQ15')
$cell.Value = 'pathogenicity'
$formatSource.Copy()
$cell.PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$cell.AddComment('To what is the entity pathogenic') | Out-Null

$cell = $ws.Range('AR15')
$cell.Value = 'ref_biomaterial'
$formatSource.Copy()
$cell.PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$cell.AddComment('Primary publication or genome report in the form of pubmed ID, DOI or URL') | Out-Null

$cell = $ws.Range('AS15')
$cell.Value = 'samp_collect_device'
$formatSource.Copy()
$cell.PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$cell.AddComment('Method or device employed for collecting sample') | Out-Null

$cell = $ws.Range('AT15')
$cell.Value = 'samp_mat_process'
$formatSource.Copy()
$cell.PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$cell.AddComment('Processing applied to the sample during or after isolation') | Out-Null

$cell = $ws.Range('AU15')
$cell.Value = 'samp_size'
$formatSource.Copy()
$cell.PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$cell.AddComment('Amount or size of sample (volume, mass or area) that was collected') | Out-Null

$cell = $ws.Range('AV15')
$cell.Value = 'samp_sort_meth'
$formatSource.Copy()
$cell.PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$cell.AddComment('method by which samples are sorted') | Out-Null

$cell = $ws.Range('AW15')
$cell.Value = 'samp_vol_we_dna_ext'
$formatSource.Copy()
$cell.PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$cell.AddComment('volume (mL) or weight (g) of sample processed for DNA extraction') | Out-Null

$cell = $ws.Range('AX15')
$cell.Value = 'source_material_id'
$formatSource.Copy()
$cell.PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$cell.AddComment('unique identifier assigned to a material sample used for extracting nucleic acids, and subsequent sequencing. The identifier can refer either to the original material collected or to any derived sub-samples.') | Out-Null

$cell = $ws.Range('AY15')
$cell.Value = 'subspecf_gen_lin'
$formatSource.Copy()
$cell.PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$cell.AddComment('Information about the genetic distinctness of the lineage (eg., biovar, serovar)') | Out-Null

$cell = $ws.Range('AZ15')
$cell.Value = 'substructure_type'
$formatSource.Copy()
$cell.PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$cell.AddComment('substructure or under building is that largely hidden section of the building which is built off the foundations to the ground floor level') | Out-Null

$cell = $ws.Range('BA15')
$cell.Value = 'surf_air_cont'
$formatSource.Copy()
$cell.PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$cell.AddComment('contaminant identified on surface') | Out-Null

$cell = $ws.Range('BB15')
$cell.Value = 'surf_humidity'
$formatSource.Copy()
$cell.PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$cell.AddComment('surfaces: water activity as a function of air and material moisture') | Out-Null

$cell = $ws.Range('BC15')
$cell.Value = 'surf_material'
$formatSource.Copy()
$cell.PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$cell.AddComment('surface materials at the point of sampling') | Out-Null

$cell = $ws.Range('BD15')
$cell.Value = 'surf_moisture'
$formatSource.Copy()
$cell.PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$cell.AddComment('water held on a surface') | Out-Null

$cell = $ws.Range('BE15')
$cell.Value = 'surf_moisture_ph'
$formatSource.Copy()
$cell.PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$cell.AddComment('pH measurement of surface') | Out-Null

$cell = $ws.Range('BF15')
$cell.Value = 'surf_temp'
$formatSource.Copy()
$cell.PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$cell.AddComment('temperature of the surface at the time of sampling') | Out-Null
